# Update the SG model on the "Apparatus" sheet of the IEEE 14-Bus workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Apparatus")

# --- Update SG model parameter values -------------------------------------
# Row 3
$ws.Range("C3").Value = 1.8
$ws.Range("E3").Value = 0.33
$ws.Range("G3").Value = 0.95
$ws.Range("K3").Value = 9.32

# Row 4
$ws.Range("K4").Value = 8.57

# Row 5
$ws.Range("C5").Value = 1.8
$ws.Range("E5").Value = 0.33
$ws.Range("G5").Value = 0.95
$ws.Range("K5").Value = 8.88

# Row 6
$ws.Range("K6").Value = 8.32

# Row 7
$ws.Range("K7").Value = 9.21

# --- Update sheet view / selection state -----------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E56").Select()

# --- Update workbook window size state --------------------------------------
$win.Width = 16200
$win.Height = 24855
